# PowerShell-style Word COM-interop script applying the proof-reading edits
# described by the commit "Proof read and edit documentations".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "Son, Dongwoo (20420487)" -- drop the spell-check split so it's one run
Replace-Text "Son, Dongwoo (20420487)" "Son, Dongwoo (20420487)"

# 2. Intro paragraph proof-reading
Replace-Text "sentences are very likely to grammatically incorrect" "sentences are very likely to be grammatically incorrect"
Replace-Text "In this assignment 2 issues of generating sentences " "In this assignment, two issues of generating sentences "
Replace-Text "national language resources will be experienced. " "national language resources will be examined. "

# 3. "We will experience this issue" -> "We will experiment this issue"
Replace-Text "We will experience this issue by specifying" "We will experiment this issue by specifying"

# 4. Harry Potter ontology sentence
Replace-Text "The Harry Potter example ontology (harrypotter0.owl) given in the class is used" "The Harry Potter example ontology (harrypotter0.owl) provided in this course is used"
Replace-Text "as mentioned in A), " "as mentioned in A) above, "

# 5. Hypothesis sentence
Replace-Text "Hypothesis for testing – there are 2 hypothesis for this research" "Hypothesis for testing – there are 2 hypothesis for this research"

# 6. NL Names paragraph
Replace-Text "It uses lexicons to generate appropriate form of strings" "It uses lexicon entries to generate appropriate form of strings"

# 6b. "Lexicon Entries" heading: split the old pPr-only run into its own blank
#     Heading2 paragraph, then merge/tidy the remaining runs.
$d.Content.Find.Execute("Many lexicon entries are added", $true, $false, $false, $false, $false, $true, 1, $false, "^pMany lexicon entries are added", 2) | Out-Null
Replace-Text "Lexicon Entries " "Lexicon Entries "
Replace-Text " (ex. Gender, tense if verb)" " (e.g., Gender, tense if verb)"

# 7. "What worked well" paragraph
Replace-Text "and (almost) grammatically correct sentences for the ontology. However" "and grammatically correct sentences for most of the texts generated based on the ontology. However"
Replace-Text "a4_partD_son-cho.pdf)" "a4_partD_son-cho.pdf)."
